# Append a new price observation row (2025-01-25, 5.83) to the price history sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowRange = $ws.Range("A88:B88")

# Force the new cells to be treated as text (matching the existing text-typed
# date/value columns further up the sheet) instead of being auto-converted to a
# date serial number / a native number by Excel's type inference.
$rowRange.NumberFormat = "@"

$ws.Range("A88").Value = "2025-01-25"
$ws.Range("B88").Value = "5.83"

# Restore the default (unstyled) cell formatting so the new row matches the
# style-less cells used by the rest of the recent price history.
$rowRange.Style = "Normal"
